$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: new metadata record (MCH340 collection) ---
$ws.Range("A2").Value = "MCH340"
$ws.Range("C2").Value = "DOCUMENTS, SA PRESS CLIPS"
$ws.Range("E2").Value = "Series"
$ws.Range("F2").Value = "1 Box"
$ws.Range("G2").Value = "LOCATION: 33I | GRAP COUNT NUMER: NONE"

# D2 and H2 are style-only (no value), matching the rest of row 2's font
foreach ($addr in @("A2","C2","D2","E2","F2","G2","H2")) {
    $c = $ws.Range($addr)
    $c.Font.Name = "Calibri"
    $c.Font.ThemeColor = 1
}

# F2 (the extent/quantity cell) is visually distinguished with its own alignment
$ws.Range("F2").HorizontalAlignment = -4108

# --- Selection / freeze-pane view state ---
$win = $excel.ActiveWindow
$ws.Range("A2").Select()
$win.FreezePanes = $true
$ws.Range("A2:H2").Select()
